$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.173420786857605
$ws.Range("B1").Value = 2.436959981918335
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.366005659103394
$ws.Range("E1").Value = 1.236374020576477
